# "Refactor Code - WEB" - adds a Transaction_Category / tran_type_query pair of
# columns to the MyAccount.xlsx test-data sheet, and moves the saved cursor /
# scroll position to reflect the new, wider layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) and new data cell (row 2) in columns G/H -------
# G1/H1 are brand-new headers; G2 already holds a value (unchanged) and H2 is
# the new data cell that goes with the new H1 header.
$ws.Range("G1").Value = "tran_type_query"
$ws.Range("H1").Value = "Transaction_Category"
$ws.Range("H2").Value = "Non Financial"

# --- Column widths -----------------------------------------------------------
# Column G shrinks slightly (new neighbour columns reduced the "best fit"
# width) and two new bestFit/customWidth columns (H, I) are introduced.
$ws.Columns("G").ColumnWidth = 185.0823333333333
$ws.Columns("H").ColumnWidth = 32.0823333333333
$ws.Columns("I").ColumnWidth = 9.1666666666667

# --- View / selection state ---------------------------------------------------
# Scroll the window so column H is the left-most visible column and move the
# active selection to M5, matching the sheet's saved view state.
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M5").Select()
